# edit.ps1 -- applies the "add 2022-Q3 data" change described by the diff.
#
# Summary of the change:
#   1. A new worksheet "2022-Q3" is inserted right after "总计" (so sheet
#      order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3,
#      2021-Q2, 2021-Q1) and populated with the fund-holdings table for
#      that quarter.
#   2. The "总计" (summary) sheet gets a new second row for 2022-Q3
#      (18 funds held, 0.44 亿元 market value), and every row below it
#      shifts down by one, with the running index in column A renumbered.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $zj)
$q3.Name = "2022-Q3"

# Seed it from the "2022-Q2" sheet's layout (header row + index column +
# styling) so column widths/styles match the other quarterly sheets,
# then overwrite every cell with the 2022-Q3 figures below.
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Range("A1:H19").Copy($q3.Range("A1:H19"))

# Numeric-looking columns that must stay text (fund codes / percentages /
# amounts with significant leading zeros or fixed decimal formatting).
$q3.Range("B2:B19").NumberFormat = "@"
$q3.Range("D2:G19").NumberFormat = "@"

$q3Data = @(
    @('167508', '安信价值发现两年定期开放混合（LOF）', '3.26', '89.75', '3.41', '0.1112', 8),
    @('006377', '广发趋势动力灵活配置混合', '1.52', '90.52', '4.69', '0.0713', 3),
    @('003029', '安信新优选灵活配置混合C', '3.61', '33.25', '1.76', '0.0635', 9),
    @('009907', '湘财长泽灵活配置混合A', '0.93', '87.32', '4.73', '0.0440', 6),
    @('002025', '广发聚盛灵活配置混合A', '3.56', '22.01', '0.79', '0.0281', 1),
    @('009766', '安信平稳双利3个月持有期混合A', '0.92', '30.99', '2.31', '0.0213', 7),
    @('005544', '银华瑞和灵活配置混合', '0.63', '89.62', '3.30', '0.0208', 8),
    @('004393', '安信企业价值优选混合', '0.54', '87.43', '3.44', '0.0186', 9),
    @('009908', '湘财长泽灵活配置混合C', '0.39', '87.32', '4.73', '0.0184', 6),
    @('012673', '华融融兴6个月定开混合A', '0.45', '36.77', '3.55', '0.0160', 1),
    @('001399', '安信鑫安得利灵活配置混合A', '0.70', '30.15', '1.67', '0.0117', 8),
    @('001400', '安信鑫安得利灵活配置混合C', '0.25', '30.15', '1.67', '0.0042', 8),
    @('009767', '安信平稳双利3个月持有期混合C', '0.17', '30.99', '2.31', '0.0039', 7),
    @('002026', '广发聚盛灵活配置混合C', '0.46', '22.01', '0.79', '0.0036', 1),
    @('750005', '安信平稳增长混合A', '0.07', '58.57', '4.19', '0.0029', 4),
    @('002035', '安信平稳增长混合C', '0.03', '58.57', '4.19', '0.0013', 4),
    @('003028', '安信新优选灵活配置混合A', '0.06', '33.25', '1.76', '0.0011', 9),
    @('012674', '华融融兴6个月定开混合C', '0.02', '36.77', '3.55', '0.0007', 1)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Range("A$r").Value = $r - 2
    $q3.Range("B$r").Value = $row[0]
    $q3.Range("C$r").Value = $row[1]
    $q3.Range("D$r").Value = $row[2]
    $q3.Range("E$r").Value = $row[3]
    $q3.Range("F$r").Value = $row[4]
    $q3.Range("G$r").Value = $row[5]
    $q3.Range("H$r").Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the 2022-Q3 row and shift
#    the existing quarters down by one row.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Shift the existing data rows (2-7) down to (3-8), carrying their
# styling along, then overwrite every row (2-8) with the correct final
# values so there is no residual floating point drift from the copy.
$zj.Range("A2:D7").Copy($zj.Range("A3:D8"))

$totalData = @(
    @("2022-Q3", 18, 0.44),
    @("2022-Q2", 32, 3.03),
    @("2022-Q1", 29, 2.4),
    @("2021-Q4", 4, 0.2),
    @("2021-Q3", 6, 0.18),
    @("2021-Q2", 3, 0.15),
    @("2021-Q1", 2, 0.01)
)

$r = 2
foreach ($row in $totalData) {
    $zj.Range("A$r").Value = $r - 2
    $zj.Range("B$r").Value = $row[0]
    $zj.Range("C$r").Value = $row[1]
    $zj.Range("D$r").Value = $row[2]
    $r = $r + 1
}
